# Append one new data row to each of the four worksheets, matching the
# rows that were appended to the source logging database.
# Column layout: A=time (datetime), B=len, C=payload, D=len2, E=checksum,
#                F..I decimal decodes.

$wb = $excel.ActiveWorkbook

$sheetRows = @(
    @{ Sheet = "ROW50-FE-LIFTER";  Row = 47; A = 45747.68463601852; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x66"; E = "0xe";  F = 400; G = 568631262647114.0 * 1000000000; H = 358; I = 14 },
    @{ Sheet = "ROW11-FE-LIFTER";  Row = 47; A = 45747.70941738426; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x01,0x66"; E = "0x14"; F = 400; G = 568631262647114.0 * 1000000000; H = 358; I = 20 },
    @{ Sheet = "ROW11-MID-LIFTER"; Row = 47; A = 45747.84657114583; B = "0x01,0x90";  C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x01,0x6e"; E = "0x19"; F = 400; G = 568631262647114.0 * 1000000000; H = 366; I = 25 }
)

foreach ($entry in $sheetRows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}

# ROW50-MID-LIFTER keeps column G as literal text (matches the rest of
# that sheet's existing rows, which store the big integer as a string).
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r2 = 49

$ws2.Cells.Item($r2, 1).Value = 45747.65561342592
$ws2.Cells.Item($r2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws2.Cells.Item($r2, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r2, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r2, 4).Value = "0x01,0x6e"
$ws2.Cells.Item($r2, 5).Value = "0x19"
$ws2.Cells.Item($r2, 6).Value = 400
$ws2.Cells.Item($r2, 7).NumberFormat = "@"
$ws2.Cells.Item($r2, 7).Value = "568631262647113771663628"
$ws2.Cells.Item($r2, 7).Style = "Normal"
$ws2.Cells.Item($r2, 8).Value = 366
$ws2.Cells.Item($r2, 9).Value = 25
